$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-08-30 Saturday" "2025-08-31 Sunday"

Replace-Text "270×5=1350" "382×9=3438"
Replace-Text "749×4=2996" "452×4=1808"
Replace-Text "824×4=3296" "529×6=3174"
Replace-Text "986×7=6902" "524×9=4716"
Replace-Text "883×5=4415" "947×6=5682"

Replace-Text "625×2=1250" "437×4=1748"
Replace-Text "596×2=1192" "909×7=6363"
Replace-Text "780×3=2340" "908×7=6356"
Replace-Text "890×7=6230" "736×5=3680"
Replace-Text "980×3=2940" "158×4=632"

Replace-Text "938×3=2814" "617×2=1234"
Replace-Text "189×9=1701" "993×8=7944"
Replace-Text "315×4=1260" "152×9=1368"
Replace-Text "125×7=875" "664×5=3320"
Replace-Text "431×3=1293" "551×9=4959"

Replace-Text "966×5=4830" "221×2=442"
Replace-Text "255×6=1530" "270×2=540"
Replace-Text "536×5=2680" "329×4=1316"
Replace-Text "939×5=4695" "307×4=1228"
Replace-Text "413×7=2891" "683×9=6147"

Replace-Text "897×3=2691" "838×6=5028"
Replace-Text "149×4=596" "606×5=3030"
Replace-Text "912×9=8208" "735×8=5880"
Replace-Text "889×2=1778" "611×9=5499"
Replace-Text "647×2=1294" "612×9=5508"
